$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture a plain/default cell style to re-apply after writes, so that
# forcing text (via a leading apostrophe) does not leave a stray
# quote-prefix / number-format style on the cell.
$defaultStyle = $ws.Range("A1").Style

$ws.Range("D2").Value = "'26.226.95"
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Value = "'  -0.79%  "
$ws.Range("E2").Style = $defaultStyle
$ws.Range("D3").Value = "'1.655.45"
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Value = "'  -1.16%  "
$ws.Range("E3").Style = $defaultStyle
$ws.Range("E4").Value = "'  -0.66%  "
$ws.Range("E4").Style = $defaultStyle
$ws.Range("D5").Value = "'219.34"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "'  -0.89%  "
$ws.Range("E5").Style = $defaultStyle
$ws.Range("D6").Value = "'0.5229"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "'  -2.38%  "
$ws.Range("E6").Style = $defaultStyle
$ws.Range("D7").Value = "'1.004"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = "'  -0.67%  "
$ws.Range("E7").Style = $defaultStyle
$ws.Range("D8").Value = "'0.2661"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "'  -0.55%  "
$ws.Range("E8").Style = $defaultStyle
$ws.Range("D9").Value = "'0.06341"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "'  -1.15%  "
$ws.Range("E9").Style = $defaultStyle
$ws.Range("E10").Value = "'  -1.93%  "
$ws.Range("E10").Style = $defaultStyle
$ws.Range("D11").Value = "'0.07773"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "'  -1.04%  "
$ws.Range("E11").Style = $defaultStyle
$ws.Range("D12").Value = "'4.549"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "'  +0.00%  "
$ws.Range("E12").Style = $defaultStyle
$ws.Range("D13").Value = "'1.677.04"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = "'  +0.05%  "
$ws.Range("E13").Style = $defaultStyle
$ws.Range("D14").Value = "'1.884.23"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "'  -1.04%  "
$ws.Range("E14").Style = $defaultStyle
$ws.Range("D15").Value = "'0.5674"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "'  +0.17%  "
$ws.Range("E15").Style = $defaultStyle
$ws.Range("D16").Value = "'0.0₅8112"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "'  -1.19%  "
$ws.Range("E16").Style = $defaultStyle
$ws.Range("D17").Value = "'65.45"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "'  -1.56%  "
$ws.Range("E17").Style = $defaultStyle
$ws.Range("D18").Value = "'26.225.56"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "'  -0.95%  "
$ws.Range("E18").Style = $defaultStyle
$ws.Range("E19").Value = "'  -0.67%  "
$ws.Range("E19").Style = $defaultStyle
$ws.Range("D20").Value = "'4.717"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "'  -0.45%  "
$ws.Range("E20").Style = $defaultStyle
$ws.Range("D21").Value = "'192.38"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "'  -3.02%  "
$ws.Range("E21").Style = $defaultStyle
$ws.Range("E22").Value = "'  -0.57%  "
$ws.Range("E22").Style = $defaultStyle
$ws.Range("D23").Value = "'6.031"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "'  -0.96%  "
$ws.Range("E23").Style = $defaultStyle
$ws.Range("E24").Value = "'  -0.65%  "
$ws.Range("E24").Style = $defaultStyle
$ws.Range("D25").Value = "'143.59"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "'  -2.20%  "
$ws.Range("E25").Style = $defaultStyle
$ws.Range("D26").Value = "'0.1200"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "'  -2.81%  "
$ws.Range("E26").Style = $defaultStyle
$ws.Range("D27").Value = "'7.269"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "'  -0.05%  "
$ws.Range("E27").Style = $defaultStyle
$ws.Range("D28").Value = "'15.98"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = "'  -1.71%  "
$ws.Range("E28").Style = $defaultStyle
$ws.Range("D29").Value = "'1.495"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "'  -0.99%  "
$ws.Range("E29").Style = $defaultStyle
$ws.Range("D30").Value = "'0.05606"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "'  -4.86%  "
$ws.Range("E30").Style = $defaultStyle
$ws.Range("D31").Value = "'1.279"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "'  -0.73%  "
$ws.Range("E31").Style = $defaultStyle
$ws.Range("D32").Value = "'3.503"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "'  -2.44%  "
$ws.Range("E32").Style = $defaultStyle
$ws.Range("D33").Value = "'3.378"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "'  +1.76%  "
$ws.Range("E33").Style = $defaultStyle
$ws.Range("D34").Value = "'1.583"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "'  -2.39%  "
$ws.Range("E34").Style = $defaultStyle
$ws.Range("D35").Value = "'2.802"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "'  -1.71%  "
$ws.Range("E35").Style = $defaultStyle
$ws.Range("D36").Value = "'0.9456"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "'  -2.70%  "
$ws.Range("E36").Style = $defaultStyle
$ws.Range("D37").Value = "'2.401"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "'  -1.59%  "
$ws.Range("E37").Style = $defaultStyle
$ws.Range("D38").Value = "'0.5773"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "'  -1.01%  "
$ws.Range("E38").Style = $defaultStyle
$ws.Range("E39").Value = "'  -1.19%  "
$ws.Range("E39").Style = $defaultStyle
$ws.Range("D40").Value = "'5.904"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("D41").Value = "'2.579"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "'  -0.71%  "
$ws.Range("E41").Style = $defaultStyle
$ws.Range("D42").Value = "'0.8470"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "'  -2.37%  "
$ws.Range("E42").Style = $defaultStyle
$ws.Range("E43").Value = "'  -0.74%  "
$ws.Range("E43").Style = $defaultStyle
$ws.Range("D44").Value = "'1.032.32"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "'  -4.45%  "
$ws.Range("E44").Style = $defaultStyle
$ws.Range("D45").Value = "'102.18"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "'  -2.13%  "
$ws.Range("E45").Style = $defaultStyle
$ws.Range("D46").Value = "'1.795.16"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "'  -1.02%  "
$ws.Range("E46").Style = $defaultStyle
$ws.Range("D47").Value = "'58.43"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "'  -0.40%  "
$ws.Range("E47").Style = $defaultStyle
$ws.Range("B48").Value = "'Frax"
$ws.Range("B48").Style = $defaultStyle
$ws.Range("C48").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("C48").Style = $defaultStyle
$ws.Range("D48").Value = "'1.001"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "'  -1.35%  "
$ws.Range("E48").Style = $defaultStyle
$ws.Range("B49").Value = "'Cronos"
$ws.Range("B49").Style = $defaultStyle
$ws.Range("C49").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C49").Style = $defaultStyle
$ws.Range("D49").Value = "'0.05314"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "'  +2.79%  "
$ws.Range("E49").Style = $defaultStyle
$ws.Range("B50").Value = "'BabyDogeCoin"
$ws.Range("B50").Style = $defaultStyle
$ws.Range("C50").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C50").Style = $defaultStyle
$ws.Range("D50").Value = "'0.0₈103"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "'  -3.82%  "
$ws.Range("E50").Style = $defaultStyle
$ws.Range("D51").Value = "'0.4353"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "'  -1.40%  "
$ws.Range("E51").Style = $defaultStyle
